$d = $word.ActiveDocument

$replacements = @(
    @("521÷8=65, 1", "787÷4=196, 3"),
    @("535÷3=178, 1", "339÷6=56, 3"),
    @("836÷3=278, 2", "745÷4=186, 1"),
    @("186÷5=37, 1", "119÷3=39, 2"),
    @("346÷4=86, 2", "784÷2=392, 0"),
    @("802÷9=89, 1", "391÷2=195, 1"),
    @("216÷9=24, 0", "179÷9=19, 8"),
    @("776÷6=129, 2", "323÷2=161, 1"),
    @("410÷4=102, 2", "256÷6=42, 4"),
    @("474÷9=52, 6", "267÷4=66, 3"),
    @("138÷2=69, 0", "260÷7=37, 1"),
    @("767÷8=95, 7", "540÷2=270, 0"),
    @("621÷2=310, 1", "678÷4=169, 2"),
    @("310÷2=155, 0", "301÷9=33, 4"),
    @("746÷9=82, 8", "369÷3=123, 0"),
    @("641÷5=128, 1", "450÷8=56, 2"),
    @("748÷9=83, 1", "308÷8=38, 4"),
    @("940÷9=104, 4", "314÷2=157, 0"),
    @("394÷7=56, 2", "120÷6=20, 0"),
    @("688÷5=137, 3", "558÷5=111, 3"),
    @("847÷9=94, 1", "592÷9=65, 7"),
    @("816÷3=272, 0", "448÷4=112, 0"),
    @("843÷6=140, 3", "577÷7=82, 3"),
    @("710÷7=101, 3", "673÷2=336, 1"),
    @("529÷7=75, 4", "889÷4=222, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
